# Adds unit-id mappings for the EasySolar-II/Multiplus-II GX VE.Direct ports
# (ttyS1 / ttyS2) to the "Unit ID mapping" sheet, and logs the change on the
# "Document versions" sheet.
#
# 273 -> 230 (ttyS1)
# 274 -> 229 (ttyS2)
# https://github.com/victronenergy/venus/issues/616

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Unit ID mapping" sheet: insert two new rows at the top of the table
#    (row 7), pushing the existing entries down, and fill in the new
#    mappings.
# ---------------------------------------------------------------------------
$wsMap = $wb.Worksheets.Item("Unit ID mapping")

$wsMap.Rows.Item(7).Insert()
$wsMap.Rows.Item(7).Insert()

$wsMap.Range("A7").Value = 230
$wsMap.Range("B7").Value = 273
$wsMap.Range("C7").Value = "EasySolar-II/Multiplus-II GX VE.Direct port (ttyS1)"

$wsMap.Range("A8").Value = 229
$wsMap.Range("B8").Value = 274
$wsMap.Range("C8").Value = "EasySolar-II/Multiplus-II GX VE.Direct port (ttyS2)"

# ---------------------------------------------------------------------------
# 2. "Document versions" sheet: append the new revision entry.
# ---------------------------------------------------------------------------
$wsVer = $wb.Worksheets.Item("Document versions")

$wsVer.Range("A69").Value = "Rev 28"
$wsVer.Range("B69").Value = "Added mapping for EasySolar-II/Multiplus-II GX VE.Direct ports"

# ---------------------------------------------------------------------------
# 3. Restore view state: "Unit ID mapping" scrolled back to the top with
#    B7 selected, and "Document versions" (the originally active tab)
#    scrolled to the top with A1 selected and re-activated last so it stays
#    the active sheet.
# ---------------------------------------------------------------------------
$wsMap.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsMap.Range("B7").Select() | Out-Null

$wsVer.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsVer.Range("A1").Select() | Out-Null
